$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 is a zero-padded text code ("001"); force text so Excel doesn't
# coerce it to the number 1, then clear the format change so no new
# cell style is introduced.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("J2").ClearFormats()

$ws.Range("M2").Value = "2020-12-18 00:00:00"
$ws.Range("N2").Value = "2017-12-31 00:00:00"

$ws.Range("O2").Value = 214805333.33
$ws.Range("P2").Value = 116.2797548748
$ws.Range("Q2").Value = 2590629118.28
$ws.Range("R2").Value = 1402.3754167328
$ws.Range("S2").Value = 229224350.34
$ws.Range("T2").Value = 124.0851465635
$ws.Range("U2").Value = 471520003.25
$ws.Range("V2").Value = 255.2461316789
$ws.Range("Y2").Value = 105282103.12
$ws.Range("Z2").Value = 56.9919608313
$ws.Range("AA2").Value = -501603581.52
$ws.Range("AB2").Value = -271.5311607923
$ws.Range("AC2").Value = 184731498.24
